$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) format, used to restore
# the style index after a NumberFormat detour for text-looking numbers.
$defaultStyle = $ws.Range("C2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.838.45"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.737.02"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.14"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +3.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5124"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2735"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06102"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.737.13"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07161"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.95"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6355"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.590"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.18"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.861.01"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.65"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006735"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  +2.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.959.76"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.258"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  +2.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.660"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.225"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.75"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.515"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.11"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.750"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.38"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  +3.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.976"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  +8.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08328"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.637"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +4.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04556"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.670"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  +2.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9831"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6166"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01597"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.927"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.001"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.99"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  -1.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3834"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7360"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  +2.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.950"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05264"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.163"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.74"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  +3.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.49"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.545"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3411"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  +1.75%  "
